# Append the 2019-04-23 ("23/4") region-fan rows to the bottom of the
# continuously-growing RegionDF sheet (rows 331-337), mirroring the
# existing A:C layout (Date Fetched | Region | Fans).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel serial date number for 2019-04-23 (matches the 43577 = 2019-04-22
# rows already on the sheet).
$serialDate = 43578

$data = @(
    @("Attica (region), Greece", 163),
    @("Central Macedonia, Greece", 114),
    @("Thessaly, Greece", 26),
    @("Western Greece, Greece", 22),
    @("Central Greece (region), Greece", 20),
    @("Eastern Macedonia and Thrace, Greece", 20),
    @("Epirus (region), Greece", 15)
)

$startRow = 331
$lastRow = $startRow - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $region = $data[$i][0]
    $fans = $data[$i][1]

    $ws.Cells.Item($row, 1).Value = $serialDate
    $ws.Cells.Item($row, 2).Value = $region
    $ws.Cells.Item($row, 3).Value = $fans

    # Reuse the date-number-format style already used by the existing
    # "Date Fetched" column instead of minting a brand-new numFmt.
    $ws.Cells.Item($lastRow, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null # xlPasteFormats
}

$excel.CutCopyMode = $false

$ws.Range("C332").Select()
